$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10; this shifts the old rows 10 and 11 down to 11 and 12.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new weekly price entry.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44841
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100107
$ws.Cells.Item(10, 8).Value = "Otros"
$ws.Cells.Item(10, 9).Value = 100107002
$ws.Cells.Item(10, 10).Value = "Chirimoya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 60
$ws.Cells.Item(10, 14).Value = 23000
$ws.Cells.Item(10, 15).Value = 24000
$ws.Cells.Item(10, 16).Value = 23500
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 19).Value = 2350
$ws.Cells.Item(10, 20).Value = 10
